$d = $word.ActiveDocument

$pairs = @(
    @("409×8=", "597×8="),
    @("325×2=", "340×7="),
    @("635×7=", "902×5="),
    @("717×7=", "227×4="),
    @("514×4=", "463×2="),
    @("182×4=", "635×3="),
    @("616×9=", "126×2="),
    @("709×4=", "267×6="),
    @("481×7=", "256×7="),
    @("513×7=", "385×3="),
    @("193×9=", "436×9="),
    @("257×4=", "702×9="),
    @("234×8=", "125×6="),
    @("773×9=", "463×5="),
    @("352×8=", "246×7="),
    @("234×3=", "662×6="),
    @("189×3=", "302×6="),
    @("510×4=", "648×3="),
    @("323×6=", "186×5="),
    @("956×7=", "985×8="),
    @("722×2=", "225×2="),
    @("304×3=", "405×8="),
    @("606×6=", "767×9="),
    @("394×2=", "902×3="),
    @("670×8=", "655×3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
